# Apply weekly refresh of Fruta/Hortaliza data (Camote - Terminal La Palmera de La Serena)
# The underlying data rows (2..11) got reshuffled into a new order; only the
# Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) columns actually change
# per row (the rest of the columns are identical across all rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per row: Row -> @{ D; J; K; L; M; P }
$data = @{
    2  = @{ D = 45005; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    3  = @{ D = 44964; J = 300; K = 20000; L = 21000; M = 20500; P = 1139 }
    4  = @{ D = 44960; J = 400; K = 19500; L = 20000; M = 19750; P = 1097 }
    5  = @{ D = 44998; J = 320; K = 17000; L = 18000; M = 17500; P = 972 }
    6  = @{ D = 44568; J = 500; K = 15000; L = 16000; M = 15500; P = 861 }
    7  = @{ D = 44957; J = 400; K = 21000; L = 22000; M = 21500; P = 1194 }
    8  = @{ D = 44547; J = 200; K = 13000; L = 14000; M = 13500; P = 750 }
    9  = @{ D = 44977; J = 400; K = 16500; L = 17000; M = 16750; P = 931 }
    10 = @{ D = 44557; J = 400; K = 13000; L = 14000; M = 13500; P = 750 }
    11 = @{ D = 44984; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
